$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E1").Value = "Spacing - 1mm"
$ws.Range("E2").Value = "Velocity of tool - 5mm/s"
$ws.Range("E3").Value = "Transition time- .25s"

$ws.Range("E3").Select()
